$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate/swap existing row contents (columns F:V) per source diff ---
# Row 2 <- old Row 4 data
$ws.Cells.Item(2,6).Value = 'Knokke'
$ws.Cells.Item(2,7).Value = 3
$ws.Cells.Item(2,8).Value = 'Cappellen'
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(2,10).Value = 1.71
$ws.Cells.Item(2,11).Value = '30/08/2023 10:16'
$ws.Cells.Item(2,12).Value = 1.65
$ws.Cells.Item(2,13).Value = '30/08/2023 19:08'
$ws.Cells.Item(2,14).Value = 3.89
$ws.Cells.Item(2,15).Value = '30/08/2023 10:16'
$ws.Cells.Item(2,16).Value = 4.43
$ws.Cells.Item(2,17).Value = '30/08/2023 19:08'
$ws.Cells.Item(2,18).Value = 4.24
$ws.Cells.Item(2,19).Value = '30/08/2023 10:16'
$ws.Cells.Item(2,20).Value = 4.17
$ws.Cells.Item(2,21).Value = '30/08/2023 19:08'
$ws.Cells.Item(2,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/knokke-cappellen/pM7B7FT7/'

# Row 4 <- old Row 6 data
$ws.Cells.Item(4,6).Value = 'Tienen'
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 'Heist'
$ws.Cells.Item(4,9).Value = 1
$ws.Cells.Item(4,10).Value = 2.92
$ws.Cells.Item(4,11).Value = '30/08/2023 09:12'
$ws.Cells.Item(4,12).Value = 2.53
$ws.Cells.Item(4,13).Value = '30/08/2023 19:24'
$ws.Cells.Item(4,14).Value = 3.36
$ws.Cells.Item(4,15).Value = '30/08/2023 09:12'
$ws.Cells.Item(4,16).Value = 3.55
$ws.Cells.Item(4,17).Value = '30/08/2023 18:06'
$ws.Cells.Item(4,18).Value = 2.23
$ws.Cells.Item(4,19).Value = '30/08/2023 09:12'
$ws.Cells.Item(4,20).Value = 2.54
$ws.Cells.Item(4,21).Value = '30/08/2023 19:24'
$ws.Cells.Item(4,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/tienen-heist/CSacVILE/'

# Row 6 <- old Row 2 data
$ws.Cells.Item(6,6).Value = 'Dessel'
$ws.Cells.Item(6,7).Value = 1
$ws.Cells.Item(6,8).Value = 'Namur'
$ws.Cells.Item(6,9).Value = 0
$ws.Cells.Item(6,10).Value = 1.98
$ws.Cells.Item(6,11).Value = '30/08/2023 09:42'
$ws.Cells.Item(6,12).Value = 1.88
$ws.Cells.Item(6,13).Value = '30/08/2023 19:27'
$ws.Cells.Item(6,14).Value = 3.61
$ws.Cells.Item(6,15).Value = '30/08/2023 09:42'
$ws.Cells.Item(6,16).Value = 4
$ws.Cells.Item(6,17).Value = '30/08/2023 19:27'
$ws.Cells.Item(6,18).Value = 3.28
$ws.Cells.Item(6,19).Value = '30/08/2023 09:42'
$ws.Cells.Item(6,20).Value = 3.45
$ws.Cells.Item(6,21).Value = '30/08/2023 19:27'
$ws.Cells.Item(6,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/dessel-namur/fik2UxyL/'

# Row 70 <- old Row 71 data
$ws.Cells.Item(70,6).Value = 'Virton'
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = 'Namur'
$ws.Cells.Item(70,9).Value = 1
$ws.Cells.Item(70,10).Value = 1.65
$ws.Cells.Item(70,11).Value = '21/10/2023 11:12'
$ws.Cells.Item(70,12).Value = 1.57
$ws.Cells.Item(70,13).Value = '21/10/2023 19:47'
$ws.Cells.Item(70,14).Value = 3.91
$ws.Cells.Item(70,15).Value = '21/10/2023 11:12'
$ws.Cells.Item(70,16).Value = 4.19
$ws.Cells.Item(70,17).Value = '21/10/2023 19:47'
$ws.Cells.Item(70,18).Value = 4.43
$ws.Cells.Item(70,19).Value = '21/10/2023 11:12'
$ws.Cells.Item(70,20).Value = 5.12
$ws.Cells.Item(70,21).Value = '21/10/2023 19:47'
$ws.Cells.Item(70,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/virton-namur/jindQod8/'

# Row 71 <- old Row 72 data
$ws.Cells.Item(71,6).Value = 'RAAL La Louviere'
$ws.Cells.Item(71,7).Value = 4
$ws.Cells.Item(71,8).Value = 'OC Charleroi'
$ws.Cells.Item(71,9).Value = 0
$ws.Cells.Item(71,10).Value = 1.36
$ws.Cells.Item(71,11).Value = '21/10/2023 11:12'
$ws.Cells.Item(71,12).Value = 1.38
$ws.Cells.Item(71,13).Value = '21/10/2023 19:56'
$ws.Cells.Item(71,14).Value = 4.86
$ws.Cells.Item(71,15).Value = '21/10/2023 11:12'
$ws.Cells.Item(71,16).Value = 4.91
$ws.Cells.Item(71,17).Value = '21/10/2023 19:56'
$ws.Cells.Item(71,18).Value = 6.78
$ws.Cells.Item(71,19).Value = '21/10/2023 11:12'
$ws.Cells.Item(71,20).Value = 7.11
$ws.Cells.Item(71,21).Value = '21/10/2023 19:56'
$ws.Cells.Item(71,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/raal-la-louviere-oc-charleroi/EZ2D12J7/'

# Row 72 <- old Row 70 data
$ws.Cells.Item(72,6).Value = 'Hoogstraten'
$ws.Cells.Item(72,7).Value = 2
$ws.Cells.Item(72,8).Value = 'Dessel'
$ws.Cells.Item(72,9).Value = 1
$ws.Cells.Item(72,10).Value = 2.15
$ws.Cells.Item(72,11).Value = '20/10/2023 07:12'
$ws.Cells.Item(72,12).Value = 2.09
$ws.Cells.Item(72,13).Value = '21/10/2023 19:57'
$ws.Cells.Item(72,14).Value = 3.36
$ws.Cells.Item(72,15).Value = '20/10/2023 07:12'
$ws.Cells.Item(72,16).Value = 3.58
$ws.Cells.Item(72,17).Value = '21/10/2023 19:39'
$ws.Cells.Item(72,18).Value = 2.85
$ws.Cells.Item(72,19).Value = '20/10/2023 07:12'
$ws.Cells.Item(72,20).Value = 3.2
$ws.Cells.Item(72,21).Value = '21/10/2023 19:39'
$ws.Cells.Item(72,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/hoogstraten-dessel/zsDX6cem/'

# Row 73 <- old Row 74 data
$ws.Cells.Item(73,6).Value = 'Thes Sport'
$ws.Cells.Item(73,7).Value = 0
$ws.Cells.Item(73,8).Value = 'Leuven B'
$ws.Cells.Item(73,9).Value = 3
$ws.Cells.Item(73,10).Value = 1.67
$ws.Cells.Item(73,11).Value = '21/10/2023 19:19'
$ws.Cells.Item(73,12).Value = 1.67
$ws.Cells.Item(73,13).Value = '21/10/2023 19:19'
$ws.Cells.Item(73,14).Value = 4.15
$ws.Cells.Item(73,15).Value = '21/10/2023 19:19'
$ws.Cells.Item(73,16).Value = 4.15
$ws.Cells.Item(73,17).Value = '21/10/2023 19:19'
$ws.Cells.Item(73,18).Value = 4.33
$ws.Cells.Item(73,19).Value = '21/10/2023 19:19'
$ws.Cells.Item(73,20).Value = 4.33
$ws.Cells.Item(73,21).Value = '21/10/2023 19:19'
$ws.Cells.Item(73,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/thes-sport-leuven/zq492r41/'

# Row 74 <- old Row 73 data
$ws.Cells.Item(74,6).Value = 'Heist'
$ws.Cells.Item(74,7).Value = 1
$ws.Cells.Item(74,8).Value = 'Lokeren-Temse'
$ws.Cells.Item(74,9).Value = 2
$ws.Cells.Item(74,10).Value = 2.74
$ws.Cells.Item(74,11).Value = '20/10/2023 07:12'
$ws.Cells.Item(74,12).Value = 3.98
$ws.Cells.Item(74,13).Value = '21/10/2023 18:33'
$ws.Cells.Item(74,14).Value = 3.29
$ws.Cells.Item(74,15).Value = '20/10/2023 07:12'
$ws.Cells.Item(74,16).Value = 3.53
$ws.Cells.Item(74,17).Value = '21/10/2023 18:37'
$ws.Cells.Item(74,18).Value = 2.29
$ws.Cells.Item(74,19).Value = '20/10/2023 07:12'
$ws.Cells.Item(74,20).Value = 1.86
$ws.Cells.Item(74,21).Value = '21/10/2023 18:36'
$ws.Cells.Item(74,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/heist-sc-lokeren-temse/Uso0P5BE/'

# Row 82 <- old Row 83 data
$ws.Cells.Item(82,6).Value = 'Virton'
$ws.Cells.Item(82,7).Value = 0
$ws.Cells.Item(82,8).Value = 'Thes Sport'
$ws.Cells.Item(82,9).Value = 1
$ws.Cells.Item(82,10).Value = 1.77
$ws.Cells.Item(82,11).Value = '27/10/2023 08:13'
$ws.Cells.Item(82,12).Value = 1.77
$ws.Cells.Item(82,13).Value = '28/10/2023 02:16'
$ws.Cells.Item(82,14).Value = 3.54
$ws.Cells.Item(82,15).Value = '27/10/2023 08:13'
$ws.Cells.Item(82,16).Value = 3.73
$ws.Cells.Item(82,17).Value = '28/10/2023 19:05'
$ws.Cells.Item(82,18).Value = 3.87
$ws.Cells.Item(82,19).Value = '27/10/2023 08:13'
$ws.Cells.Item(82,20).Value = 4.2
$ws.Cells.Item(82,21).Value = '28/10/2023 19:05'
$ws.Cells.Item(82,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/virton-thes-sport/vugc1gPO/'

# Row 83 <- old Row 82 data
$ws.Cells.Item(83,6).Value = 'OC Charleroi'
$ws.Cells.Item(83,7).Value = 2
$ws.Cells.Item(83,8).Value = 'Hoogstraten'
$ws.Cells.Item(83,9).Value = 1
$ws.Cells.Item(83,10).Value = 2.47
$ws.Cells.Item(83,11).Value = '28/10/2023 13:13'
$ws.Cells.Item(83,12).Value = 2.1
$ws.Cells.Item(83,13).Value = '28/10/2023 19:45'
$ws.Cells.Item(83,14).Value = 3.31
$ws.Cells.Item(83,15).Value = '28/10/2023 13:13'
$ws.Cells.Item(83,16).Value = 3.37
$ws.Cells.Item(83,17).Value = '28/10/2023 19:45'
$ws.Cells.Item(83,18).Value = 2.62
$ws.Cells.Item(83,19).Value = '28/10/2023 13:13'
$ws.Cells.Item(83,20).Value = 3.36
$ws.Cells.Item(83,21).Value = '28/10/2023 19:45'
$ws.Cells.Item(83,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/oc-charleroi-hoogstraten/hrsuQx2P/'

# Row 84 <- old Row 85 data
$ws.Cells.Item(84,6).Value = 'Namur'
$ws.Cells.Item(84,7).Value = 2
$ws.Cells.Item(84,8).Value = 'URSL Vise'
$ws.Cells.Item(84,9).Value = 0
$ws.Cells.Item(84,10).Value = 2.24
$ws.Cells.Item(84,11).Value = '29/10/2023 12:42'
$ws.Cells.Item(84,12).Value = 2.43
$ws.Cells.Item(84,13).Value = '29/10/2023 14:05'
$ws.Cells.Item(84,14).Value = 3.38
$ws.Cells.Item(84,15).Value = '29/10/2023 12:42'
$ws.Cells.Item(84,16).Value = 3.37
$ws.Cells.Item(84,17).Value = '29/10/2023 14:05'
$ws.Cells.Item(84,18).Value = 2.89
$ws.Cells.Item(84,19).Value = '29/10/2023 12:42'
$ws.Cells.Item(84,20).Value = 2.75
$ws.Cells.Item(84,21).Value = '29/10/2023 14:05'
$ws.Cells.Item(84,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/namur-ursl-vise/rRvNHE95/'

# Row 85 <- old Row 84 data
$ws.Cells.Item(85,6).Value = 'Antwerp B'
$ws.Cells.Item(85,7).Value = 0
$ws.Cells.Item(85,8).Value = 'Lokeren-Temse'
$ws.Cells.Item(85,9).Value = 1
$ws.Cells.Item(85,10).Value = 5.73
$ws.Cells.Item(85,11).Value = '28/10/2023 03:42'
$ws.Cells.Item(85,12).Value = 6.25
$ws.Cells.Item(85,13).Value = '29/10/2023 14:20'
$ws.Cells.Item(85,14).Value = 4.29
$ws.Cells.Item(85,15).Value = '28/10/2023 03:42'
$ws.Cells.Item(85,16).Value = 4.31
$ws.Cells.Item(85,17).Value = '29/10/2023 14:20'
$ws.Cells.Item(85,18).Value = 1.41
$ws.Cells.Item(85,19).Value = '28/10/2023 03:42'
$ws.Cells.Item(85,20).Value = 1.47
$ws.Cells.Item(85,21).Value = '29/10/2023 14:20'
$ws.Cells.Item(85,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/antwerp-sc-lokeren-temse/fiT6LdAt/'

# --- Append new rows 89-92 ---
# First, copy row 88 (A:V) into each new row to replicate cell styles (s="1" for col A, s="2" for col E)
$ws.Range("A88:V88").Copy($ws.Range("A89:V89"))
$ws.Range("A88:V88").Copy($ws.Range("A90:V90"))
$ws.Range("A88:V88").Copy($ws.Range("A91:V91"))
$ws.Range("A88:V88").Copy($ws.Range("A92:V92"))

# Now overwrite values for the new rows
# Row 89
$ws.Cells.Item(89,1).Value = 88
$ws.Cells.Item(89,2).Value = 'belgium'
$ws.Cells.Item(89,3).Value = 'national-division-1'
$ws.Cells.Item(89,4).Value = '2023-2024'
$ws.Cells.Item(89,5).Value = 45234.83333333334
$ws.Cells.Item(89,6).Value = 'RAAL La Louviere'
$ws.Cells.Item(89,7).Value = 3
$ws.Cells.Item(89,8).Value = 'Dessel'
$ws.Cells.Item(89,9).Value = 0
$ws.Cells.Item(89,10).Value = 1.22
$ws.Cells.Item(89,11).Value = '03/11/2023 08:12'
$ws.Cells.Item(89,12).Value = 1.14
$ws.Cells.Item(89,13).Value = '04/11/2023 19:55'
$ws.Cells.Item(89,14).Value = 5.71
$ws.Cells.Item(89,15).Value = '03/11/2023 08:12'
$ws.Cells.Item(89,16).Value = 7.82
$ws.Cells.Item(89,17).Value = '04/11/2023 19:55'
$ws.Cells.Item(89,18).Value = 8.39
$ws.Cells.Item(89,19).Value = '03/11/2023 08:12'
$ws.Cells.Item(89,20).Value = 15.17
$ws.Cells.Item(89,21).Value = '04/11/2023 19:55'
$ws.Cells.Item(89,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/raal-la-louviere-dessel/4tTVFhvI/'

# Row 90
$ws.Cells.Item(90,1).Value = 89
$ws.Cells.Item(90,2).Value = 'belgium'
$ws.Cells.Item(90,3).Value = 'national-division-1'
$ws.Cells.Item(90,4).Value = '2023-2024'
$ws.Cells.Item(90,5).Value = 45234.83333333334
$ws.Cells.Item(90,6).Value = 'Lokeren-Temse'
$ws.Cells.Item(90,7).Value = 2
$ws.Cells.Item(90,8).Value = 'Cappellen'
$ws.Cells.Item(90,9).Value = 1
$ws.Cells.Item(90,10).Value = 1.2
$ws.Cells.Item(90,11).Value = '03/11/2023 08:12'
$ws.Cells.Item(90,12).Value = 1.15
$ws.Cells.Item(90,13).Value = '04/11/2023 19:55'
$ws.Cells.Item(90,14).Value = 5.89
$ws.Cells.Item(90,15).Value = '03/11/2023 08:12'
$ws.Cells.Item(90,16).Value = 7.64
$ws.Cells.Item(90,17).Value = '04/11/2023 19:55'
$ws.Cells.Item(90,18).Value = 8.77
$ws.Cells.Item(90,19).Value = '03/11/2023 08:12'
$ws.Cells.Item(90,20).Value = 14.76
$ws.Cells.Item(90,21).Value = '04/11/2023 19:59'
$ws.Cells.Item(90,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/sc-lokeren-temse-cappellen/jX1lYZgn/'

# Row 91
$ws.Cells.Item(91,1).Value = 90
$ws.Cells.Item(91,2).Value = 'belgium'
$ws.Cells.Item(91,3).Value = 'national-division-1'
$ws.Cells.Item(91,4).Value = '2023-2024'
$ws.Cells.Item(91,5).Value = 45234.83333333334
$ws.Cells.Item(91,6).Value = 'Charleroi B'
$ws.Cells.Item(91,7).Value = 1
$ws.Cells.Item(91,8).Value = 'Gent B'
$ws.Cells.Item(91,9).Value = 2
$ws.Cells.Item(91,10).Value = 2.56
$ws.Cells.Item(91,11).Value = '04/11/2023 16:09'
$ws.Cells.Item(91,12).Value = 3.14
$ws.Cells.Item(91,13).Value = '04/11/2023 18:23'
$ws.Cells.Item(91,14).Value = 3.38
$ws.Cells.Item(91,15).Value = '04/11/2023 16:09'
$ws.Cells.Item(91,16).Value = 3.36
$ws.Cells.Item(91,17).Value = '04/11/2023 18:13'
$ws.Cells.Item(91,18).Value = 2.56
$ws.Cells.Item(91,19).Value = '04/11/2023 16:09'
$ws.Cells.Item(91,20).Value = 2.2
$ws.Cells.Item(91,21).Value = '04/11/2023 18:23'
$ws.Cells.Item(91,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/charleroi-gent/C6WwEW8U/'

# Row 92
$ws.Cells.Item(92,1).Value = 91
$ws.Cells.Item(92,2).Value = 'belgium'
$ws.Cells.Item(92,3).Value = 'national-division-1'
$ws.Cells.Item(92,4).Value = '2023-2024'
$ws.Cells.Item(92,5).Value = 45234.85416666666
$ws.Cells.Item(92,6).Value = 'Thes Sport'
$ws.Cells.Item(92,7).Value = 0
$ws.Cells.Item(92,8).Value = 'OC Charleroi'
$ws.Cells.Item(92,9).Value = 3
$ws.Cells.Item(92,10).Value = 2.49
$ws.Cells.Item(92,11).Value = '04/11/2023 16:09'
$ws.Cells.Item(92,12).Value = 2.49
$ws.Cells.Item(92,13).Value = '04/11/2023 16:09'
$ws.Cells.Item(92,14).Value = 3.37
$ws.Cells.Item(92,15).Value = '04/11/2023 18:34'
$ws.Cells.Item(92,16).Value = 3.37
$ws.Cells.Item(92,17).Value = '04/11/2023 18:34'
$ws.Cells.Item(92,18).Value = 2.64
$ws.Cells.Item(92,19).Value = '04/11/2023 16:09'
$ws.Cells.Item(92,20).Value = 2.64
$ws.Cells.Item(92,21).Value = '04/11/2023 16:09'
$ws.Cells.Item(92,22).Value = 'https://www.betexplorer.com/football/belgium/national-division-1/thes-sport-oc-charleroi/QkURGYOB/'

